# Applies the "Updated cryptos list" data refresh described in the commit:
# - Refreshes Price (column D) and Volume(1h) (column E) percentages for most rows.
# - Swaps rows 30 and 31 (Kaspa <-> EthereumClassic), including Coin name, Link, Price and Volume.
#
# Every value is written as literal text (a leading apostrophe forces Excel to treat
# numeric-looking strings like "1.50" or "19.24" as text instead of coercing them to a
# number), and the cell style is reset to "Normal" afterwards so the quote-prefix marker
# doesn't leave a stray style behind (matching the original workbook, where these cells
# carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Ref,
        [string]$Value
    )
    $cell = $ws.Range($Ref)
    $cell.Value = "'" + $Value
    $cell.Style = "Normal"
}


Set-CellText "D2" "37.361.64"
Set-CellText "E2" "  +2.33%  "
Set-CellText "D3" "2.066.28"
Set-CellText "E3" "  +3.61%  "
Set-CellText "E4" "  +0.03%  "
Set-CellText "D5" "235.72"
Set-CellText "E5" "  +0.71%  "
Set-CellText "D6" "0.614"
Set-CellText "E6" "  +2.61%  "
Set-CellText "D7" "58.47"
Set-CellText "E7" "  +6.75%  "
Set-CellText "E8" "  +0.02%  "
Set-CellText "E9" "  +3.38%  "
Set-CellText "D10" "58.97"
Set-CellText "E11" "  +2.12%  "
Set-CellText "E12" "  +2.72%  "
Set-CellText "D13" "2.371.78"
Set-CellText "E13" "  +3.65%  "
Set-CellText "D14" "14.59"
Set-CellText "E14" "  +3.14%  "
Set-CellText "E15" "  +5.23%  "
Set-CellText "E16" "  +2.68%  "
Set-CellText "D17" "5.19"
Set-CellText "E17" "  +2.50%  "
Set-CellText "D18" "2.084.38"
Set-CellText "E18" "  +4.52%  "
Set-CellText "D19" "37.405.13"
Set-CellText "E19" "  +2.57%  "
Set-CellText "E20" "  +16.16%  "
Set-CellText "D21" "70.35"
Set-CellText "D22" "0.0₃0815"
Set-CellText "E22" "  +1.36%  "
Set-CellText "D23" "227.47"
Set-CellText "E23" "  +2.52%  "
Set-CellText "E24" "  +0.14%  "
Set-CellText "E25" "  +2.46%  "
Set-CellText "E26" "  +0.93%  "
Set-CellText "D27" "165.45"
Set-CellText "E27" "  +2.24%  "
Set-CellText "D28" "1.50"
Set-CellText "E28" "  +13.74%  "
Set-CellText "D29" "8.89"
Set-CellText "E29" "  +2.39%  "
Set-CellText "B30" "EthereumClassic"
Set-CellText "C30" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText "D30" "19.24"
Set-CellText "E30" "  +2.07%  "
Set-CellText "B31" "Kaspa"
Set-CellText "C31" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText "D31" "0.127"
Set-CellText "E31" "  +0.12%  "
Set-CellText "E32" "  +2.55%  "
Set-CellText "E33" "  +2.92%  "
Set-CellText "E34" "  +3.41%  "
Set-CellText "D35" "2.54"
Set-CellText "E35" "  +8.21%  "
Set-CellText "D36" "4.56"
Set-CellText "E36" "  +7.52%  "
Set-CellText "D38" "1.00"
Set-CellText "E38" "  -0.01%  "
Set-CellText "E39" "  +1.41%  "
Set-CellText "E40" "  +5.03%  "
Set-CellText "D41" "0.0979"
Set-CellText "E41" "  +2.96%  "
Set-CellText "E42" "  -1.19%  "
Set-CellText "D43" "4.46"
Set-CellText "E43" "  +24.56%  "
Set-CellText "D44" "1.458.25"
Set-CellText "E44" "  +0.60%  "
Set-CellText "D45" "95.63"
Set-CellText "E45" "  +7.51%  "
Set-CellText "E46" "  +6.46%  "
Set-CellText "E47" "  +4.59%  "
Set-CellText "D48" "15.88"
Set-CellText "E48" "  +4.62%  "
Set-CellText "E49" "  +4.19%  "
Set-CellText "D50" "7.28"
Set-CellText "E50" "  +6.40%  "
Set-CellText "E51" "  +2.13%  "
